# Daily attendance processing - 2025-11-14 22:47:25
# Reorders the "Recorded By" (column G) comma-separated list so that any
# list containing a "system"/"System" entry has its parts reversed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notlike "*,*") { continue }

    $parts = $val -split ","
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $hasSystem = $false
    foreach ($p in $trimmed) {
        if ($p.ToLower() -eq "system") { $hasSystem = $true }
    }

    if (-not $hasSystem) { continue }

    $reversed = $trimmed[($trimmed.Count - 1)..0]
    $newVal = [string]::Join(", ", $reversed)

    $cell.Value = $newVal
}
